$d = $word.ActiveDocument

# 1) "Gravitationskraft am Fg" -> "Gravitationskraft Fg"
$d.Content.Find.Execute("eine Gravitationskraft am ", $true, $false, $false, $false, $false, $true, 1, $false, "eine Gravitationskraft ", 2) | Out-Null

# 2) "ein Kräftegleichgewicht M = M einstellt" -> "ein Gleichgewicht der Drehmomente einstellt"
$d.Content.Find.Execute("Kräftegleichgewicht M = M einstellt", $true, $false, $false, $false, $false, $true, 1, $false, "Gleichgewicht der Drehmomente einstellt", 2) | Out-Null

# 3) "ausgedrückt wurde." -> "ausgedrückt worden."
$d.Content.Find.Execute("ausgedrückt wurde.", $true, $false, $false, $false, $false, $true, 1, $false, "ausgedrückt worden.", 2) | Out-Null

# 4) Helmholtz coil paragraph expansion
$d.Content.Find.Execute("Flussdichte … erzeugt,", $true, $false, $false, $false, $false, $true, 1, $false, "Flussdichte … in der Mitte erzeugt,", 2) | Out-Null

$d.Content.Find.Execute("Superposition zweier Spulen darstellen lässt. Das Feld hat die Eigenschaft, dass es im Inneren der Spule nahezu homogen ist.", $true, $false, $false, $false, $false, $true, 1, $false, "Superposition zweier Spulen darstellen lässt. Dabei ist µ_0 die magnetische Feldkonstante, I der die Spulen durchfließende Strom, R ihr Radius und x ihr halber Abstand. Das Feld hat die Eigenschaft, dass es im Inneren auf der Symmetrieachse der Spulen nahezu homogen ist.", 2) | Out-Null

# 5) "Fassung, in der die Kugel" -> "Fassung, in die die Kugel"
$d.Content.Find.Execute("Fassung, in der die Kugel", $true, $false, $false, $false, $false, $true, 1, $false, "Fassung, in die die Kugel", 2) | Out-Null
